# "changed binary cost files to €/g" — the filament cost in economical_params!B2
# was stored as a per-kg price (45.5) and is converted to a per-gram price
# (0.076 EUR/g). The Scaling sheet's B2/C2 formulas reference B2 and simply
# recalculate to the new, smaller values.

$wb = $excel.ActiveWorkbook

$wsParams  = $wb.Worksheets.Item("economical_params")
$wsScaling = $wb.Worksheets.Item("Scaling")

# Core data edit: total_cost_filament, €/kg -> €/g
$wsParams.Range("B2").Value = 0.076

# Force a recalculation so the dependent formulas on Scaling pick up the
# new value immediately.
$excel.Calculate()

# UI state: the author ended up with economical_params active (selection on
# B22) instead of Scaling (which previously was the active tab, selection D7).
$wsParams.Activate()
$wsParams.Range("B22").Select()

Write-Output "Updated economical_params!B2 to 0.076 EUR/g and refreshed Scaling formulas"
